# Update the lattice-multiplication practice table: each cell holds a
# 5-line "problem" (title "A x B", the two multiplier digits, a dashed
# rule, and two empty-box rows seeded with the digits of A). This swaps
# every cell's problem for a new A x B pair while keeping the template.

function Set-LatticeCell {
    param($table, $row, $col, $title, $mid, $d1, $d2)

    $nl = [char]11   # <w:br/> line break inside a Word Range.Text assignment
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $title + $nl + $mid + $nl + "  ----" + $nl + $d1 + $nl + $d2
}

$d = $word.ActiveDocument
$t = $d.Tables(1)

Set-LatticeCell $t 1 1 "56 x 45" "  4    5" "5|    |" "6|    |"
Set-LatticeCell $t 1 2 "35 x 89" "  8    9" "3|    |" "5|    |"
Set-LatticeCell $t 1 3 "50 x 82" "  8    2" "5|    |" "0|    |"
Set-LatticeCell $t 2 1 "11 x 85" "  8    5" "1|    |" "1|    |"
Set-LatticeCell $t 2 2 "46 x 46" "  4    6" "4|    |" "6|    |"
Set-LatticeCell $t 2 3 "92 x 43" "  4    3" "9|    |" "2|    |"
Set-LatticeCell $t 3 1 "33 x 52" "  5    2" "3|    |" "3|    |"
Set-LatticeCell $t 3 2 "79 x 27" "  2    7" "7|    |" "9|    |"
Set-LatticeCell $t 3 3 "26 x 82" "  8    2" "2|    |" "6|    |"
Set-LatticeCell $t 4 1 "28 x 40" "  4    0" "2|    |" "8|    |"
Set-LatticeCell $t 4 2 "65 x 48" "  4    8" "6|    |" "5|    |"
Set-LatticeCell $t 4 3 "95 x 74" "  7    4" "9|    |" "5|    |"
Set-LatticeCell $t 5 1 "90 x 67" "  6    7" "9|    |" "0|    |"
Set-LatticeCell $t 5 2 "57 x 35" "  3    5" "5|    |" "7|    |"
Set-LatticeCell $t 5 3 "69 x 13" "  1    3" "6|    |" "9|    |"
